# Add the new "2022-Q1" sheet between "2021-Q4" and "总计", and populate
# it with the quarterly fund-holding detail, then prepend a matching
# summary row onto the "总计" sheet.
#
# NOTE: worksheet object handles returned by this COM layer behave as
# *positional* references - once a structural operation (Add/Move/Delete)
# changes sheet ordering, any previously-captured handle can start
# pointing at a different sheet. To stay safe, all structural changes are
# performed first, and every worksheet reference used afterwards is
# (re-)fetched fresh, by name, only once ordering is final.

$wb = $excel.ActiveWorkbook

# --- create + position the new sheet (structural changes go first) -------
$wsTotalBeforeAdd = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q1"

$wsTotalForMove = $wb.Worksheets.Item("总计")
$wsQ1ForMove = $wb.Worksheets.Item("2022-Q1")
$wsQ1ForMove.Move($wsTotalForMove, $null)

# Ordering is now final: "2021-Q4", "2022-Q1", "总计".
# Re-fetch every sheet we still need, fresh, by name.
$wsQ4 = $wb.Worksheets.Item("2021-Q4")
$wsQ1 = $wb.Worksheets.Item("2022-Q1")
$wsTotal = $wb.Worksheets.Item("总计")

# Touch a cell so the freshly created sheet is "realized" before any
# cross-sheet range copy is attempted against it.
$wsQ1.Cells.Item(1, 1).Value = "x"

# Header row (B1:H1) - identical text/format to the "2021-Q4" sheet, so
# simply copy it over (brings along the bold/bordered/centered style).
$wsQ4.Range("B1:H1").Copy($wsQ1.Range("B1:H1"))

# Column A (row index numbers 0..12) - copy format+values from the
# "2021-Q4" sheet's A2:A14 (same 0-based index pattern, 13 rows).
$wsQ4.Range("A2:A14").Copy($wsQ1.Range("A2:A14"))

# Remove the realization placeholder now that real content anchors the sheet.
$wsQ1.Cells.Item(1, 1).ClearContents()

# --- fund rows -------------------------------------------------------------
# Columns: B=fund code, C=fund name, D=fund size, E=stock position,
# F=position ratio, G=held value, H=rank.
# B, D, E, F, G are numeric-looking text -> write with a leading
# apostrophe so Excel keeps them as text instead of converting to numbers.
$rows = @(
    @("010054", "万家健康产业混合A", "8.13", "86.63", "3.89", "0.3163", 4),
    @("000780", "鹏华医疗保健股票", "7.67", "82.80", "3.15", "0.2416", 8),
    @("519125", "浦银安盛消费升级混合A", "4.64", "89.16", "4.86", "0.2255", 7),
    @("010055", "万家健康产业混合C", "3.36", "86.63", "3.89", "0.1307", 4),
    @("519176", "浦银安盛消费升级混合C", "2.33", "89.16", "4.86", "0.1132", 7),
    @("970032", "东海证券海睿进取灵活配置混合型集合资产管理计划A", "1.57", "78.62", "4.20", "0.0659", 2),
    @("970033", "东海证券海睿进取灵活配置混合型集合资产管理计划B", "1.18", "78.62", "4.20", "0.0496", 2),
    @("160921", "大成多策略混合(LOF)", "1.13", "79.19", "3.78", "0.0427", 10),
    @("008037", "兴银先锋成长混合A", "0.41", "79.32", "2.26", "0.0093", 4),
    @("008038", "兴银先锋成长混合C", "0.17", "79.32", "2.26", "0.0038", 4),
    @("009649", "嘉实精选平衡混合A", "0.06", "67.70", "3.81", "0.0023", 7),
    @("001474", "兴银丰盈灵活配置混合", "0.08", "83.45", "2.62", "0.0021", 5),
    @("009650", "嘉实精选平衡混合C", "0.01", "67.70", "3.81", "0.0004", 7)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]
    $wsQ1.Cells.Item($r, 2).Value = "'" + $data[0]
    $wsQ1.Cells.Item($r, 3).Value = $data[1]
    $wsQ1.Cells.Item($r, 4).Value = "'" + $data[2]
    $wsQ1.Cells.Item($r, 5).Value = "'" + $data[3]
    $wsQ1.Cells.Item($r, 6).Value = "'" + $data[4]
    $wsQ1.Cells.Item($r, 7).Value = "'" + $data[5]
    $wsQ1.Cells.Item($r, 8).Value = $data[6]
}

# --- update the "总计" (totals) sheet --------------------------------------
# Shift the existing 2021-Q4 summary row down (copy format+value), then
# write the new 2022-Q1 summary row into the now-vacant row 2.
$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3:D3"))
$wsTotal.Cells.Item(3, 1).Value = 1

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q1"
$wsTotal.Cells.Item(2, 3).Value = 13
$wsTotal.Cells.Item(2, 4).Value = 1.2
